# Update "想去人数" (want-to-go count) figures on sheets "展览" and "全部类型"
# to reflect freshly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 2745
$wsExpo.Range("F6").Value = 1924
$wsExpo.Range("F8").Value = 122
$wsExpo.Range("F9").Value = 974
$wsExpo.Range("F11").Value = 16

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 2745
$wsAll.Range("F6").Value = 1924
$wsAll.Range("F9").Value = 122
$wsAll.Range("F10").Value = 974
$wsAll.Range("F12").Value = 16
